$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9145299145299145
$ws.Range("C2").Value = 0.9145299145299145
$ws.Range("D2").Value = 0.9145299145299145

$ws.Range("B3").Value = 0.9180327868852459
$ws.Range("C3").Value = 0.9180327868852459
$ws.Range("D3").Value = 0.9180327868852459

$ws.Range("B4").Value = 0.9163179916317992
$ws.Range("C4").Value = 0.9163179916317992
$ws.Range("D4").Value = 0.9163179916317992
$ws.Range("E4").Value = 0.9163179916317992

$ws.Range("B5").Value = 0.9162813507075802
$ws.Range("C5").Value = 0.9162813507075802
$ws.Range("D5").Value = 0.9162813507075802

$ws.Range("B6").Value = 0.9163179916317992
$ws.Range("C6").Value = 0.9163179916317992
$ws.Range("D6").Value = 0.9163179916317992
